$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing revised values (rows 434, 435, 451, 457, 458)
$updates = @{
    434 = 1182578000000
    435 = 1180145000000
    451 = 1297962000000
    457 = 1375618000000
    458 = 1365405000000
}

foreach ($row in $updates.Keys) {
    $val = $updates[$row]
    $ws.Range("C$row").Value = $val
    $ws.Range("D$row").Value = $val
    $ws.Range("E$row").Value = $val
    $ws.Range("F$row").Value = $val
}

# Append new rows 460-462 (copy formatting from the last existing data row, A459)
$newRows = @(
    @{ Row = 460; Date = 44986.45833333334; Value = 1393379000000 },
    @{ Row = 461; Date = 45017.45833333334; Value = 1391394000000 },
    @{ Row = 462; Date = 45047.41666666666; Value = 1394414000000 }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Range("A459").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $nr.Date
    $ws.Cells.Item($r, 2).Value = "ECONOMICS:MAM2"
    $ws.Cells.Item($r, 3).Value = $nr.Value
    $ws.Cells.Item($r, 4).Value = $nr.Value
    $ws.Cells.Item($r, 5).Value = $nr.Value
    $ws.Cells.Item($r, 6).Value = $nr.Value
    $ws.Cells.Item($r, 7).Value = 0
}
